$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 161.4
$ws.Range("I12").Value = 161.4
$ws.Range("K12").Value = 161.4
$ws.Range("M12").Value = 8.599999999999994
$ws.Range("H112").Value = 2656.182
$ws.Range("I112").Value = 2837.8
$ws.Range("J112").Value = 2504.8333
$ws.Range("K112").Value = 8513.400000000001
$ws.Range("L112").Value = 7514.499899999999
$ws.Range("M112").Value = -7405.400000000001
$ws.Range("N112").Value = -9730.499899999999
$ws.Range("H125").Value = 2797.8
$ws.Range("I125").Value = 1996.3334
$ws.Range("K125").Value = 17967.0006
$ws.Range("M125").Value = -15507.0006
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 223560.16
$ws.Range("I32").Value = 242308.4
$ws.Range("K32").Value = 242308.4
$ws.Range("M32").Value = -242021.4
$ws.Range("H45").Value = 3044.75
$ws.Range("I45").Value = 1975.6
$ws.Range("J45").Value = 4826.6665
$ws.Range("K45").Value = 1975.6
$ws.Range("L45").Value = 4826.6665
$ws.Range("M45").Value = -1598.6
$ws.Range("N45").Value = -5580.6665
$ws.Range("H61").Value = 3230
$ws.Range("I61").Value = 2686.1538
$ws.Range("K61").Value = 2686.1538
$ws.Range("M61").Value = -2474.1538
$ws.Range("H110").Value = 2509.1667
$ws.Range("I110").Value = 5055
$ws.Range("K110").Value = 5055
$ws.Range("M110").Value = -3010
$ws.Range("H122").Value = 2513.513
$ws.Range("I122").Value = 2052.6897
$ws.Range("K122").Value = 6158.0691
$ws.Range("M122").Value = -3708.0691
$ws.Range("H132").Value = 6773
$ws.Range("I132").Value = 5499.625
$ws.Range("J132").Value = 7904.8887
$ws.Range("K132").Value = 16498.875
$ws.Range("L132").Value = 23714.6661
$ws.Range("M132").Value = -13968.875
$ws.Range("N132").Value = -28774.6661
$ws.Range("H136").Value = 3230
$ws.Range("I136").Value = 2686.1538
$ws.Range("K136").Value = 8058.4614
$ws.Range("M136").Value = -5508.4614
$ws.Range("H138").Value = 78960.25
$ws.Range("J138").Value = 78960.25
$ws.Range("L138").Value = 78960.25
$ws.Range("N138").Value = -89240.25
$ws.Range("H139").Value = 90510.71000000001
$ws.Range("J139").Value = 90510.71000000001
$ws.Range("L139").Value = 90510.71000000001
$ws.Range("N139").Value = -100790.71
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 69247.5
$ws.Range("J13").Value = 69247.5
$ws.Range("L13").Value = 69247.5
$ws.Range("N13").Value = -69583.5
$ws.Range("H20").Value = 60943.156
$ws.Range("I20").Value = 91257
$ws.Range("J20").Value = 8976.571
$ws.Range("K20").Value = 91257
$ws.Range("L20").Value = 8976.571
$ws.Range("M20").Value = -91010
$ws.Range("N20").Value = -9470.571
$ws.Range("H37").Value = 1834.6
$ws.Range("I37").Value = 2018.75
$ws.Range("J37").Value = 1098
$ws.Range("K37").Value = 2018.75
$ws.Range("L37").Value = 1098
$ws.Range("M37").Value = -1881.75
$ws.Range("N37").Value = -1372
$ws.Range("H81").Value = 31194.75
$ws.Range("J81").Value = 34222.57
$ws.Range("L81").Value = 34222.57
$ws.Range("N81").Value = -36344.57
$ws.Range("H82").Value = 30229.125
$ws.Range("I82").Value = 15270.75
$ws.Range("J82").Value = 45187.5
$ws.Range("K82").Value = 15270.75
$ws.Range("L82").Value = 45187.5
$ws.Range("M82").Value = -14887.75
$ws.Range("N82").Value = -45953.5
$ws.Range("H84").Value = 31194.75
$ws.Range("J84").Value = 34222.57
$ws.Range("L84").Value = 102667.71
$ws.Range("N84").Value = -113275.71
$ws.Range("H85").Value = 30229.125
$ws.Range("I85").Value = 15270.75
$ws.Range("J85").Value = 45187.5
$ws.Range("K85").Value = 15270.75
$ws.Range("L85").Value = 45187.5
$ws.Range("M85").Value = -13944.75
$ws.Range("N85").Value = -47839.5
$ws.Range("H86").Value = 3589.5186
$ws.Range("J86").Value = 3070.875
$ws.Range("L86").Value = 3070.875
$ws.Range("N86").Value = -5316.875
$ws.Range("H89").Value = 3589.5186
$ws.Range("J89").Value = 3070.875
$ws.Range("L89").Value = 15354.375
$ws.Range("N89").Value = -26586.375
$ws.Range("H134").Value = 7460.8
$ws.Range("I134").Value = 7636.5713
$ws.Range("K134").Value = 22909.7139
$ws.Range("M134").Value = -20374.7139
$ws.Range("H135").Value = 64674.2
$ws.Range("J135").Value = 64674.2
$ws.Range("L135").Value = 64674.2
$ws.Range("N135").Value = -74814.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 35894.445
$ws.Range("I6").Value = 30381.25
$ws.Range("K6").Value = 30381.25
$ws.Range("M6").Value = -30268.25
$ws.Range("H31").Value = 3730
$ws.Range("I31").Value = 3730
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3730
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3435
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 3730
$ws.Range("I34").Value = 3730
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3730
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3528
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 4576.769
$ws.Range("I58").Value = 3901.1765
$ws.Range("K58").Value = 3901.1765
$ws.Range("M58").Value = -3698.1765
$ws.Range("H134").Value = 3321.7273
$ws.Range("I134").Value = 3047.4546
$ws.Range("J134").Value = 3596
$ws.Range("K134").Value = 9142.363799999999
$ws.Range("L134").Value = 10788
$ws.Range("M134").Value = -6607.363799999999
$ws.Range("N134").Value = -15858
$ws.Range("H136").Value = 4576.769
$ws.Range("I136").Value = 3901.1765
$ws.Range("K136").Value = 11703.5295
$ws.Range("M136").Value = -9153.529500000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 60003110
$ws.Range("I55").Value = 168000480
$ws.Range("K55").Value = 504001440
$ws.Range("M55").Value = -504001263
$ws.Range("H68").Value = 2960
$ws.Range("J68").Value = 3714.2856
$ws.Range("L68").Value = 11142.8568
$ws.Range("N68").Value = -12764.8568
$ws.Range("H71").Value = 2960
$ws.Range("J71").Value = 3714.2856
$ws.Range("L71").Value = 33428.5704
$ws.Range("N71").Value = -41540.5704
$ws.Range("H119").Value = 13699.75
$ws.Range("I119").Value = 5799
$ws.Range("K119").Value = 17397
$ws.Range("M119").Value = -12559
$ws.Range("H127").Value = 7474.3335
$ws.Range("J127").Value = 7474.3335
$ws.Range("L127").Value = 22423.0005
$ws.Range("N127").Value = -32343.0005
$ws.Range("H139").Value = 7578.9473
$ws.Range("I139").Value = 3800.5
$ws.Range("K139").Value = 11401.5
$ws.Range("M139").Value = -6261.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H80").Value = 2853.7334
$ws.Range("I80").Value = 2595.6
$ws.Range("J80").Value = 3370
$ws.Range("K80").Value = 2595.6
$ws.Range("L80").Value = 3370
$ws.Range("M80").Value = -1597.6
$ws.Range("N80").Value = -5366
$ws.Range("H83").Value = 2853.7334
$ws.Range("I83").Value = 2595.6
$ws.Range("J83").Value = 3370
$ws.Range("K83").Value = 12978
$ws.Range("L83").Value = 16850
$ws.Range("M83").Value = -7986
$ws.Range("N83").Value = -26834
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2590.087
$ws.Range("I22").Value = 2112.5334
$ws.Range("K22").Value = 2112.5334
$ws.Range("M22").Value = -1817.5334
$ws.Range("H27").Value = 2590.087
$ws.Range("I27").Value = 2112.5334
$ws.Range("K27").Value = 2112.5334
$ws.Range("M27").Value = -2005.5334
$ws.Range("H55").Value = 1302.75
$ws.Range("I55").Value = 1132.8
$ws.Range("K55").Value = 1132.8
$ws.Range("M55").Value = -959.8
$ws.Range("H61").Value = 7596.1924
$ws.Range("I61").Value = 6522
$ws.Range("J61").Value = 15831.667
$ws.Range("K61").Value = 6522
$ws.Range("L61").Value = 15831.667
$ws.Range("M61").Value = -6320
$ws.Range("N61").Value = -16235.667
$ws.Range("H99").Value = 42631.332
$ws.Range("I99").Value = 42631.332
$ws.Range("K99").Value = 42631.332
$ws.Range("M99").Value = -39636.332
$ws.Range("H113").Value = 7596.1924
$ws.Range("I113").Value = 6522
$ws.Range("J113").Value = 15831.667
$ws.Range("K113").Value = 6522
$ws.Range("L113").Value = 15831.667
$ws.Range("M113").Value = -4352
$ws.Range("N113").Value = -20171.667
$ws.Range("H132").Value = 3635.375
$ws.Range("I132").Value = 3083.3333
$ws.Range("K132").Value = 9249.999899999999
$ws.Range("M132").Value = -6719.999899999999
